$wb = $excel.ActiveWorkbook

# --- Senvion sheet: replace the placeholder part number with the real one ---
$ws1 = $wb.Worksheets.Item("Senvion")
$ws1.Range("B2").Value = "REPWR-PRT-0001"
$ws1.Columns("B").ColumnWidth = 13.166666666666666

# --- Nordex sheet: add the missing header row ---
$ws2 = $wb.Worksheets.Item("Nordex")
$ws2.Range("A1").Value = "Nordex Part No"
$ws2.Range("B1").Value = "Mpulse Part No."
$ws2.Columns("A").ColumnWidth = 13.833333333333332

$ws2.Activate() | Out-Null
$ws2.Range("B1").Select() | Out-Null
